$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5958322286605835
$ws.Range("B1").Value = 0.6447725892066956
$ws.Range("D1").Value = 1.514372110366821
$ws.Range("E1").Value = 0.9047685265541077
